# daily auto push: 2025-10-10 18:38 UTC
# Append the new day's record as row 91 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 91

# Columns A (date) and B (weekday) are plain text in this sheet, not real
# dates, so force a text number format before assigning the value (Excel
# would otherwise auto-detect "2025/10/11" as a date serial number).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/11"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "土"

$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 201

$wb.Save()
